{"js": "// Replace \"20 sekunder\" with \"10 sekunder\" in the NFR paragraph\n// (\"Hovedscenariet gennemf\u00f8res p\u00e5 h\u00f8jest 20 sekunder i 95% af tilf\u00e6ldene.\")\nconst results = context.document.body.search(\"h\u00f8jest 20 sekunder\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"h\u00f8jest 10 sekunder\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace \"20 sekunder\" with \"10 sekunder\" in the NFR paragraph\n# (\"Hovedscenariet gennemf\u00f8res p\u00e5 h\u00f8jest 20 sekunder i 95% af tilf\u00e6ldene.\")\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"h\u00f8jest 20 sekunder\"\n$find.Replacement.Text = \"h\u00f8jest 10 sekunder\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
